$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.072.35"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "3.136.09"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'534.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.27%  "

$ws.Range("D6").Value = "'139.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.74%  "

$ws.Range("D9").Value = "'7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("E10").Value = "  +2.07%  "

$ws.Range("E11").Value = "  +4.17%  "

$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").Value = "3.676.33"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "'25.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("E15").Value = "  +5.42%  "

$ws.Range("D16").Value = "58.133.28"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("D17").Value = "'6.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.71%  "

$ws.Range("D18").Value = "3.133.06"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").Value = "'12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.51%  "

$ws.Range("E20").Value = "  +3.84%  "

$ws.Range("D21").Value = "'375.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.09%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("D24").Value = "'70.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("E25").Value = "  +2.75%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "0.0₃0885"
$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").Value = "'7.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.69%  "

$ws.Range("D30").Value = "'6.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.76%  "

$ws.Range("D31").Value = "'1.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("D32").Value = "'21.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.19%  "

$ws.Range("D33").Value = "'5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.26%  "

$ws.Range("E34").Value = "  +3.14%  "

$ws.Range("D35").Value = "'161.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("E36").Value = "  +3.27%  "

$ws.Range("D37").Value = "'1.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.34%  "

$ws.Range("D38").Value = "'25.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("E39").Value = "  +6.53%  "

$ws.Range("D40").Value = "2.624.03"
$ws.Range("E40").Value = "  +9.19%  "

$ws.Range("D41").Value = "'0.0673"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("E42").Value = "  +4.14%  "

$ws.Range("D43").Value = "'38.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.02%  "

$ws.Range("D44").Value = "'0.700"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").Value = "'0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "'6.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.22%  "

$ws.Range("D48").Value = "'0.982"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("D49").Value = "'0.0989"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.64%  "

$ws.Range("D50").Value = "'20.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.81%  "

$ws.Range("D51").Value = "'0.751"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.06%  "
